$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (D) and volume-change (E) values per row.
# D-column writes go through a Text-format round-trip so Excel's automatic
# General-format type inference doesn't coerce numeric-looking strings
# (e.g. "0.999", "0.0260") into actual numbers; the style is reset back to
# Normal immediately after so no new formatting is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.968.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.35%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.765.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.69%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.81%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.614'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.40%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.765.64'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.43%  '

$ws.Range("E11").Value = '  +5.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.395'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.10%  '

$ws.Range("E13").Value = '  +2.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.252.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.941.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.34%  '

$ws.Range("E17").Value = '  +8.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.764.45'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '363.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.94%  '

$ws.Range("E22").Value = '  +2.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.537'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.80%  '

$ws.Range("E26").Value = '  +6.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0919'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +20.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '175.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.65%  '

$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.69'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.79%  '

$ws.Range("E37").Value = '  +10.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.78%  '

$ws.Range("E39").Value = '  +12.12%  '

$ws.Range("E40").Value = '  +6.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '339.56'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +11.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.11%  '

$ws.Range("E46").Value = '  +6.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.650'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0260'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '137.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.68%  '

$ws.Range("E50").Value = '  +2.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.13%  '
